# Replace the working set of sequences (rows 2-33) with a freshly
# resampled batch of image/word/category cues. Column A (index) is
# untouched; columns B (random draw), C (image), D (word) and E
# (category) are rewritten row by row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 103, "dog/dog018.jpg", "stoßen", "dog"),
    @(3, 82, "house/house017.jpg", "streifen", "house"),
    @(4, 80, "dog/dog025.jpg", "proben", "dog"),
    @(5, 75, "dog/dog016.jpg", "hassen", "dog"),
    @(6, 100, "dog/dog010.jpg", "kennen", "dog"),
    @(7, 53, "dog/dog015.jpg", "meinen", "dog"),
    @(8, 62, "dog/dog008.jpg", "betteln", "dog"),
    @(9, 30, "house/house023.jpg", "kriegen", "house"),
    @(10, 44, "dog/dog006.jpg", "nullen", "dog"),
    @(11, 93, "house/house005.jpg", "duschen", "house"),
    @(12, 87, "house/house027.jpg", "quellen", "house"),
    @(13, 66, "house/house013.jpg", "bergen", "house"),
    @(14, 86, "house/house016.jpg", "heißen", "house"),
    @(15, 81, "dog/dog013.jpg", "fließen", "dog"),
    @(16, 122, "house/house000.jpg", "bauen", "house"),
    @(17, 84, "house/house030.jpg", "zögern", "house"),
    @(18, 4, "house/house022.jpg", "spüren", "house"),
    @(19, 123, "dog/dog017.jpg", "hören", "dog"),
    @(20, 79, "house/house028.jpg", "deuten", "house"),
    @(21, 105, "house/house009.jpg", "holen", "house"),
    @(22, 13, "house/house024.jpg", "münzen", "house"),
    @(23, 12, "dog/dog007.jpg", "lernen", "dog"),
    @(24, 114, "dog/dog024.jpg", "kranken", "dog"),
    @(25, 14, "dog/dog029.jpg", "passen", "dog"),
    @(26, 126, "dog/dog001.jpg", "achten", "dog"),
    @(27, 49, "dog/dog026.jpg", "herrschen", "dog"),
    @(28, 117, "house/house003.jpg", "süßen", "house"),
    @(29, 45, "house/house018.jpg", "rechnen", "house"),
    @(30, 51, "house/house031.jpg", "binden", "house"),
    @(31, 55, "dog/dog031.jpg", "grenzen", "dog"),
    @(32, 109, "house/house021.jpg", "heben", "house"),
    @(33, 25, "dog/dog002.jpg", "piepen", "dog")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
}
